# Update NATMI TPM values for the Rarres2-Ccrl2 ligand-receptor pair sheet (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5310463333333334
$ws.Range("H2").Value = 1.593139
$ws.Range("I2").Value = 0.01208304439515528
$ws.Range("J2").Value = 0.01241653163470256
$ws.Range("M2").Value = 3.373155666666667
$ws.Range("N2").Value = 10.119467
$ws.Range("O2").Value = 0.03263881027773249
$ws.Range("P2").Value = 0.03275129872373959
$ws.Range("Q2").Value = 1.791301948545889
$ws.Range("R2").Value = 16.121717536913
$ws.Range("S2").Value = 0.0003943761935908923
$ws.Range("T2").Value = 0.0004066575366809064

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5310463333333334
$ws.Range("H3").Value = 1.593139
$ws.Range("I3").Value = 0.01208304439515528
$ws.Range("J3").Value = 0.01241653163470256
$ws.Range("O3").Value = 0.007166957156089505
$ws.Range("P3").Value = 0.007191657807437632
$ws.Range("Q3").Value = 0.3933410626675556
$ws.Range("R3").Value = 3.540069564008
$ws.Range("S3").Value = [double]"8.659866149520535E-05"
$ws.Range("T3").Value = [double]"8.929544667200503E-05"

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5310463333333334
$ws.Range("H4").Value = 1.593139
$ws.Range("I4").Value = 0.01208304439515528
$ws.Range("J4").Value = 0.01241653163470256
$ws.Range("M4").Value = 47.16807033333333
$ws.Range("N4").Value = 141.504211
$ws.Range("O4").Value = 0.4564004306085713
$ws.Range("P4").Value = 0.4579733977222395
$ws.Range("Q4").Value = 25.04843080092545
$ws.Range("R4").Value = 225.435877208329
$ws.Range("S4").Value = 0.005514706665011356
$ws.Range("T4").Value = 0.005686441180670405

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5310463333333334
$ws.Range("H5").Value = 1.593139
$ws.Range("I5").Value = 0.01208304439515528
$ws.Range("J5").Value = 0.01241653163470256
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 1.064885
$ws.Range("N5").Value = 2.12977
$ws.Range("O5").Value = 0.01030387652312217
$ws.Range("P5").Value = 0.006892925633618734
$ws.Range("Q5").Value = 0.5655032746716668
$ws.Range("R5").Value = 3.393019648030001
$ws.Range("S5").Value = 0.0001245021974710834
$ws.Range("T5").Value = [double]"8.558622918547922E-05"

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5310463333333334
$ws.Range("H6").Value = 1.593139
$ws.Range("I6").Value = 0.01208304439515528
$ws.Range("J6").Value = 0.01241653163470256
$ws.Range("M6").Value = 51.00119533333333
$ws.Range("N6").Value = 153.003586
$ws.Range("O6").Value = 0.4934899254344846
$ws.Range("P6").Value = 0.4951907201129645
$ws.Range("Q6").Value = 27.08399777738378
$ws.Range("R6").Value = 243.755979996454
$ws.Range("S6").Value = 0.005962860677586748
$ws.Range("T6").Value = 0.006148551241493766

# Row 7
$ws.Range("I7").Value = 0.8986139302737502
$ws.Range("J7").Value = 0.9234153188332328
$ws.Range("M7").Value = 3.373155666666667
$ws.Range("N7").Value = 10.119467
$ws.Range("O7").Value = 0.03263881027773249
$ws.Range("P7").Value = 0.03275129872373959
$ws.Range("Q7").Value = 133.2188173483212
$ws.Range("R7").Value = 1198.969356134891
$ws.Range("S7").Value = 0.02932968958313247
$ws.Range("T7").Value = 0.03024305095318445

# Row 8
$ws.Range("I8").Value = 0.8986139302737502
$ws.Range("J8").Value = 0.9234153188332328
$ws.Range("O8").Value = 0.007166957156089505
$ws.Range("P8").Value = 0.007191657807437632
$ws.Range("S8").Value = 0.006440327538137169
$ws.Range("T8").Value = 0.006640886987194529

# Row 9
$ws.Range("I9").Value = 0.8986139302737502
$ws.Range("J9").Value = 0.9234153188332328
$ws.Range("M9").Value = 47.16807033333333
$ws.Range("N9").Value = 141.504211
$ws.Range("O9").Value = 0.4564004306085713
$ws.Range("P9").Value = 0.4579733977222395
$ws.Range("Q9").Value = 1862.847483886978
$ws.Range("R9").Value = 16765.6273549828
$ws.Range("S9").Value = 0.4101277847278003
$ws.Range("T9").Value = 0.4228996510748207

# Row 10
$ws.Range("I10").Value = 0.8986139302737502
$ws.Range("J10").Value = 0.9234153188332328
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 1.064885
$ws.Range("N10").Value = 2.12977
$ws.Range("O10").Value = 0.01030387652312217
$ws.Range("P10").Value = 0.006892925633618734
$ws.Range("Q10").Value = 42.05638112520167
$ws.Range("R10").Value = 252.33828675121
$ws.Range("S10").Value = 0.009259206979498235
$ws.Range("T10").Value = 0.006365033121661806

# Row 11
$ws.Range("I11").Value = 0.8986139302737502
$ws.Range("J11").Value = 0.9234153188332328
$ws.Range("M11").Value = 51.00119533333333
$ws.Range("N11").Value = 153.003586
$ws.Range("O11").Value = 0.4934899254344846
$ws.Range("P11").Value = 0.4951907201129645
$ws.Range("Q11").Value = 2014.232249284686
$ws.Range("R11").Value = 18128.09024356218
$ws.Range("S11").Value = 0.4434569214451821
$ws.Range("T11").Value = 0.4572666966963713

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.541247
$ws.Range("H12").Value = 7.082494000000001
$ws.Range("I12").Value = 0.08057497440313205
$ws.Range("J12").Value = 0.0551992078554295
$ws.Range("M12").Value = 3.373155666666667
$ws.Range("N12").Value = 10.119467
$ws.Range("O12").Value = 0.03263881027773249
$ws.Range("P12").Value = 0.03275129872373959
$ws.Range("Q12").Value = 11.94517738511633
$ws.Range("R12").Value = 71.671064310698
$ws.Range("S12").Value = 0.002629871302676979
$ws.Range("T12").Value = 0.001807845745786965

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.541247
$ws.Range("H13").Value = 7.082494000000001
$ws.Range("I13").Value = 0.08057497440313205
$ws.Range("J13").Value = 0.0551992078554295
$ws.Range("O13").Value = 0.007166957156089505
$ws.Range("P13").Value = 0.007191657807437632
$ws.Range("Q13").Value = 2.622968601261333
$ws.Range("R13").Value = 15.737811607568
$ws.Range("S13").Value = 0.0005774773894002559
$ws.Range("T13").Value = 0.0003969738141378722

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3.541247
$ws.Range("H14").Value = 7.082494000000001
$ws.Range("I14").Value = 0.08057497440313205
$ws.Range("J14").Value = 0.0551992078554295
$ws.Range("M14").Value = 47.16807033333333
$ws.Range("N14").Value = 141.504211
$ws.Range("O14").Value = 0.4564004306085713
$ws.Range("P14").Value = 0.4579733977222395
$ws.Range("Q14").Value = 167.0337875637057
$ws.Range("R14").Value = 1002.202725382234
$ws.Range("S14").Value = 0.03677445301386408
$ws.Range("T14").Value = 0.02527976877312718

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.541247
$ws.Range("H15").Value = 7.082494000000001
$ws.Range("I15").Value = 0.08057497440313205
$ws.Range("J15").Value = 0.0551992078554295
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 1.064885
$ws.Range("N15").Value = 2.12977
$ws.Range("O15").Value = 0.01030387652312217
$ws.Range("P15").Value = 0.006892925633618734
$ws.Range("Q15").Value = 3.771020811595001
$ws.Range("R15").Value = 15.08408324638
$ws.Range("S15").Value = 0.0008302345871036019
$ws.Range("T15").Value = 0.0003804840347821386

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 3.541247
$ws.Range("H16").Value = 7.082494000000001
$ws.Range("I16").Value = 0.08057497440313205
$ws.Range("J16").Value = 0.0551992078554295
$ws.Range("M16").Value = 51.00119533333333
$ws.Range("N16").Value = 153.003586
$ws.Range("O16").Value = 0.4934899254344846
$ws.Range("P16").Value = 0.4951907201129645
$ws.Range("Q16").Value = 180.6078299705807
$ws.Range("R16").Value = 1083.646979823484
$ws.Range("S16").Value = 0.03976293811008714
$ws.Range("T16").Value = 0.02733413548759534

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.3835953333333333
$ws.Range("H17").Value = 1.150786
$ws.Range("I17").Value = 0.008728050927962449
$ws.Range("J17").Value = 0.008968941676635136
$ws.Range("M17").Value = 3.373155666666667
$ws.Range("N17").Value = 10.119467
$ws.Range("O17").Value = 0.03263881027773249
$ws.Range("P17").Value = 0.03275129872373959
$ws.Range("Q17").Value = 1.293926772340222
$ws.Range("R17").Value = 11.645340951062
$ws.Range("S17").Value = 0.0002848731983321534
$ws.Range("T17").Value = 0.0002937444880872752

# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.3835953333333333
$ws.Range("H18").Value = 1.150786
$ws.Range("I18").Value = 0.008728050927962449
$ws.Range("J18").Value = 0.008968941676635136
$ws.Range("O18").Value = 0.007166957156089505
$ws.Range("P18").Value = 0.007191657807437632
$ws.Range("Q18").Value = 0.2841254831768888
$ws.Range("R18").Value = 2.557129348592
$ws.Range("S18").Value = [double]"6.255356705687412E-05"
$ws.Range("T18").Value = [double]"6.450155943322584E-05"

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.3835953333333333
$ws.Range("H19").Value = 1.150786
$ws.Range("I19").Value = 0.008728050927962449
$ws.Range("J19").Value = 0.008968941676635136
$ws.Range("M19").Value = 47.16807033333333
$ws.Range("N19").Value = 141.504211
$ws.Range("O19").Value = 0.4564004306085713
$ws.Range("P19").Value = 0.4579733977222395
$ws.Range("Q19").Value = 18.09345166220511
$ws.Range("R19").Value = 162.841064959846
$ws.Range("S19").Value = 0.003983486201895602
$ws.Range("T19").Value = 0.004107536693621193

# Row 20
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.3835953333333333
$ws.Range("H20").Value = 1.150786
$ws.Range("I20").Value = 0.008728050927962449
$ws.Range("J20").Value = 0.008968941676635136
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 1.064885
$ws.Range("N20").Value = 2.12977
$ws.Range("O20").Value = 0.01030387652312217
$ws.Range("P20").Value = 0.006892925633618734
$ws.Range("Q20").Value = 0.4084849165366667
$ws.Range("R20").Value = 2.45090949922
$ws.Range("S20").Value = [double]"8.993275904924692E-05"
$ws.Range("T20").Value = [double]"6.182224798930972E-05"

# Row 21
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.3835953333333333
$ws.Range("H21").Value = 1.150786
$ws.Range("I21").Value = 0.008728050927962449
$ws.Range("J21").Value = 0.008968941676635136
$ws.Range("M21").Value = 51.00119533333333
$ws.Range("N21").Value = 153.003586
$ws.Range("O21").Value = 0.4934899254344846
$ws.Range("P21").Value = 0.4951907201129645
$ws.Range("Q21").Value = 19.56382052428844
$ws.Range("R21").Value = 176.074384718596
$ws.Range("S21").Value = 0.004307205201628573
$ws.Range("T21").Value = 0.004441336687504132
